$d = $word.ActiveDocument
$d.Content.Find.Execute("886÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "117÷2=", 2) | Out-Null
$d.Content.Find.Execute("669÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "268÷4=", 2) | Out-Null
$d.Content.Find.Execute("635÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "785÷3=", 2) | Out-Null
$d.Content.Find.Execute("978÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "224÷2=", 2) | Out-Null
$d.Content.Find.Execute("123÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "255÷7=", 2) | Out-Null
$d.Content.Find.Execute("121÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "709÷8=", 2) | Out-Null
$d.Content.Find.Execute("524÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "440÷5=", 2) | Out-Null
$d.Content.Find.Execute("174÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "293÷3=", 2) | Out-Null
$d.Content.Find.Execute("586÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "467÷7=", 2) | Out-Null
$d.Content.Find.Execute("518÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "337÷3=", 2) | Out-Null
$d.Content.Find.Execute("684÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "314÷3=", 2) | Out-Null
$d.Content.Find.Execute("691÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "848÷5=", 2) | Out-Null
$d.Content.Find.Execute("172÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "975÷9=", 2) | Out-Null
$d.Content.Find.Execute("797÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "429÷7=", 2) | Out-Null
$d.Content.Find.Execute("301÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "996÷8=", 2) | Out-Null
$d.Content.Find.Execute("374÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "597÷4=", 2) | Out-Null
$d.Content.Find.Execute("365÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "470÷4=", 2) | Out-Null
$d.Content.Find.Execute("506÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "598÷4=", 2) | Out-Null
$d.Content.Find.Execute("387÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "675÷4=", 2) | Out-Null
$d.Content.Find.Execute("390÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "798÷7=", 2) | Out-Null
$d.Content.Find.Execute("610÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "902÷5=", 2) | Out-Null
$d.Content.Find.Execute("933÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "747÷7=", 2) | Out-Null
$d.Content.Find.Execute("211÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "190÷4=", 2) | Out-Null
$d.Content.Find.Execute("569÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷4=", 2) | Out-Null
$d.Content.Find.Execute("770÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "424÷3=", 2) | Out-Null
